# Update Name of Algo
# Apply updated values to column D for the RandomForest result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = -7.024099999999994
$ws.Range("D21").Value = -8.856599999999997
$ws.Range("D23").Value = -7.252899999999998
$ws.Range("D25").Value = -8.426499999999995
$ws.Range("D53").Value = -6.379499999999996
$ws.Range("D57").Value = -8.474400000000001
$ws.Range("D59").Value = -8.2584
$ws.Range("D69").Value = -7.010499999999999
$ws.Range("D79").Value = -6.168099999999999
$ws.Range("D83").Value = -8.812199999999999
$ws.Range("D93").Value = -6.498299999999996
